# Weekly update for "Hortaliza, Vega Monumental Concepción - Espinaca".
# Three new daily price records are inserted into the historical series
# (rows keep shifting down as in the source feed); nothing else changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New record #1 -> becomes row 39 -------------------------------------
$ws.Rows("39").Insert()
$ws.Range("A39").Value = 11
$ws.Range("B39").Value = "Vega Monumental Concepción"
$ws.Range("C39").Value = "Bíobío"
$ws.Range("D39").Value = 44874
$ws.Range("E39").Value = 8
$ws.Range("F39").Value = 100112012
$ws.Range("G39").Value = "Espinaca"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 130
$ws.Range("K39").Value = 14000
$ws.Range("L39").Value = 15000
$ws.Range("M39").Value = 14615
$ws.Range("N39").Value = "$/cuna 10 kilos"
$ws.Range("O39").Value = "Provincia de Chacabuco"
$ws.Range("P39").Value = 1462
$ws.Range("Q39").Value = 10
$ws.Range("R39").Value = "Hortaliza"

# --- New record #2 -> becomes row 47 -------------------------------------
$ws.Rows("47").Insert()
$ws.Range("A47").Value = 11
$ws.Range("B47").Value = "Vega Monumental Concepción"
$ws.Range("C47").Value = "Bíobío"
$ws.Range("D47").Value = 44897
$ws.Range("E47").Value = 8
$ws.Range("F47").Value = 100112012
$ws.Range("G47").Value = "Espinaca"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 100
$ws.Range("K47").Value = 13000
$ws.Range("L47").Value = 14000
$ws.Range("M47").Value = 13500
$ws.Range("N47").Value = "$/cuna 10 kilos"
$ws.Range("O47").Value = "Región Metropolitana"
$ws.Range("P47").Value = 1350
$ws.Range("Q47").Value = 10
$ws.Range("R47").Value = "Hortaliza"

# --- New record #3 -> becomes row 49 -------------------------------------
$ws.Rows("49").Insert()
$ws.Range("A49").Value = 11
$ws.Range("B49").Value = "Vega Monumental Concepción"
$ws.Range("C49").Value = "Bíobío"
$ws.Range("D49").Value = 44616
$ws.Range("E49").Value = 8
$ws.Range("F49").Value = 100112012
$ws.Range("G49").Value = "Espinaca"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 60
$ws.Range("K49").Value = 14000
$ws.Range("L49").Value = 15000
$ws.Range("M49").Value = 14500
$ws.Range("N49").Value = "$/cuna 10 kilos"
$ws.Range("O49").Value = "Región Metropolitana"
$ws.Range("P49").Value = 1450
$ws.Range("Q49").Value = 10
$ws.Range("R49").Value = "Hortaliza"
